$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet (Through 2022-09-01 -> Through 2022-09-03)
$ws.Name = "Through 2022-09-03"

# Update header label in I1 (shared string "2022 (through 09-01)" -> "2022 (through 09-03)")
$ws.Range("I1").Value = "2022 (through 09-03)"

# Update data values for September 2022 (row 10) and Total (row 14)
$ws.Range("I10").Value = 19
$ws.Range("I14").Value = 1158
